# colombia: individuals + households
# Rename the existing "colombia" sheet to "colombia_hh" and add a
# duplicate sheet named "colombia_ind" right after it, which becomes
# the active sheet.

$wb = $excel.ActiveWorkbook

$colombia = $wb.Worksheets.Item("colombia")
$colombia.Name = "colombia_hh"

# Duplicate the household sheet to create the individuals sheet,
# placing the copy immediately after it.
$colombia.Copy($null, $colombia)

$newSheet = $wb.Worksheets.Item($colombia.Index + 1)
$newSheet.Name = "colombia_ind"

# Make the newly created sheet the active tab.
$newSheet.Activate()
